# "Object cat added to InOut, InOutAddition now fetches teams that are in Object"
#
# The "Импорт" sheet no longer needs a dedicated "Супервайзер*" column — it is
# removed, and the two trailing columns ("Питает*", "Длина*") shift left to
# take its place (C and D). On the "Супервайзеры" sheet, the name
# "Исфандиёр" is swapped out for "Мурод" (which already exists further down
# the list).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Импорт")
$ws2 = $wb.Worksheets.Item("Супервайзеры")

# --- "Импорт" sheet: drop the "Супервайзер*" header, shift D->C and E->D ---
$colDValue = $ws1.Cells.Item(1, 4).Value2
$colEValue = $ws1.Cells.Item(1, 5).Value2

$ws1.Cells.Item(1, 3).Value2 = $colDValue
$ws1.Cells.Item(1, 4).Value2 = $colEValue
$ws1.Cells.Item(1, 5).Clear()

# Update the active selection to reflect the new "Супервайзер" data column (C)
$ws1.Activate() | Out-Null
$ws1.Range("C2:C1048576").Select() | Out-Null

# --- "Супервайзеры" sheet: replace "Исфандиёр" with "Мурод" ---
$ws2.Cells.Item(2, 1).Value2 = "Мурод"
